$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E4").Value = 16
$ws.Range("E6").Value = 38
$ws.Range("E8").Value = 33
$ws.Range("E10").Value = 21
$ws.Range("H10").Value = 3
$ws.Range("E12").Value = 28
$ws.Range("H14").Value = 4

$ws.Range("H14").Select()
